# Update the NATMI Il7-Il7r LR-pairs sheet with new TPM-based values.
# The new data collapses the previous 6 data rows into 4 (Resolving-Mac
# replaces the old "FAPs"/target-cluster duplication), and all numeric
# columns are recalculated against the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then columns E..T
$rows = @(
    @("ECs",        "Il7", "Il7r", "Resolving-Mac", 3, 1,                  0.5282253333333333, 1.584676,            0.5375365413017625, 0.5375365413017625, 3, 1, 22.569632, 67.708896, 1, 1, 11.92185138641067,  107.296662477696,   0.5375365413017625, 0.5375365413017625),
    @("FAPs",       "Il7", "Il7r", "Resolving-Mac", 2, 0.6666666666666666, 0.2859449999999999, 0.8578349999999999, 0.2909854499642813, 0.2909854499642813, 3, 1, 22.569632, 67.708896, 1, 1, 6.453673422239999,  58.08306080015999,  0.2909854499642813, 0.2909854499642813),
    @("MuSCs",      "Il7", "Il7r", "Resolving-Mac", 1, 0.3333333333333333, 0.1276103333333333, 0.382831,            0.1298597641682559, 0.1298597641682559, 3, 1, 22.569632, 67.708896, 1, 1, 2.880118262730667,  25.921064364576,    0.1298597641682559, 0.1298597641682559),
    @("Resolving-Mac", "Il7", "Il7r", "Resolving-Mac", 1, 0.3333333333333333, 0.04089733333333333, 0.122692,        0.0416182445657004, 0.0416182445657004, 3, 1, 22.569632, 67.708896, 1, 1, 0.9230377631146667, 8.307339868031999,  0.0416182445657004, 0.0416182445657004)
)

# Remove all old data rows (2 through the end of the used range, which
# includes the two rows -- old rows 6 & 7 -- that are being dropped
# entirely) before rewriting rows 2-5 with the refreshed TPM numbers.
$usedRange = $ws.UsedRange
$lastRow = [Math]::Max($usedRange.Rows.Count, $rows.Count + 1)
$ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 20)).Clear()

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $row[$i]
    }
    $r++
}

$wb.Save()
